$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued Price cells (contain multiple dots - never numeric in Excel)
$ws.Range("D2").Value = '29.937.40'
$ws.Range("D3").Value = '1.890.99'
$ws.Range("D13").Value = '1.891.80'
$ws.Range("D16").Value = '29.922.98'
$ws.Range("D22").Value = '2.139.40'
$ws.Range("D50").Value = '2.039.00'

# Numeric-looking Price cells must be forced to text so they keep their original
# string formatting (trailing zeros, decimal digits) instead of becoming real numbers.
$numericPriceCells = @{
    'D4' = '1.001'
    'D5' = '0.8210'
    'D6' = '241.47'
    'D7' = '1.001'
    'D8' = '0.3244'
    'D9' = '26.42'
    'D10' = '0.07022'
    'D11' = '0.08035'
    'D12' = '0.7454'
    'D14' = '5.200'
    'D15' = '92.14'
    'D17' = '14.02'
    'D18' = '5.881'
    'D19' = '244.66'
    'D20' = '0.000007752'
    'D21' = '1.001'
    'D23' = '1.001'
    'D24' = '6.920'
    'D25' = '0.1551'
    'D27' = '9.190'
    'D31' = '1.518'
    'D32' = '4.266'
    'D33' = '0.05629'
    'D34' = '4.064'
    'D35' = '1.269'
    'D36' = '0.7281'
    'D37' = '2.716'
    'D38' = '0.01911'
    'D39' = '2.780'
    'D40' = '0.4420'
    'D41' = '71.81'
    'D42' = '5.958'
    'D43' = '0.8432'
    'D44' = '1.0000'
    'D45' = '1.872'
    'D46' = '7.568'
    'D47' = '100.54'
    'D48' = '9.744'
    'D49' = '990.82'
    'D51' = '35.90'
}
foreach ($ref in $numericPriceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericPriceCells[$ref]
    $cell.Style = "Normal"
}

# Volume(1h) percentage cells (plain text with padding spaces, e.g. "  +0.35%  ")
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +5.56%  '
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +5.91%  '
$ws.Range("E9").Value = '  +3.54%  '
$ws.Range("E10").Value = '  +2.42%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +17.78%  '
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  +7.26%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  -3.53%  '
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  +7.17%  '
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("E51").Value = '  -0.53%  '
